$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '75.857.89'
$ws.Range("E2").Value = '  +0.90%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.904.90'
$ws.Range("E3").Value = '  +1.94%  '
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '197.03'
$ws.Range("E5").Value = '  +3.76%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '595.22'
$ws.Range("E6").Value = '  -2.16%  '
$ws.Range("E7").Value = '  +0.00%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.550'
$ws.Range("E8").Value = '  -3.29%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.193'
$ws.Range("E9").Value = '  -1.85%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '2.904.63'
$ws.Range("E10").Value = '  +2.00%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.419'
$ws.Range("E11").Value = '  +11.55%  '
$ws.Range("E12").Value = '  -1.56%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.87'
$ws.Range("E13").Value = '  -1.38%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '3.428.84'
$ws.Range("E14").Value = '  +1.65%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '75.715.53'
$ws.Range("E15").Value = '  +0.78%  '
$ws.Range("B16").Value = 'ShibaInu'
$ws.Range("C16").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.0000189'
$ws.Range("E16").Value = '  -1.71%  '
$ws.Range("B17").Value = 'Avalanche'
$ws.Range("C17").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '27.26'
$ws.Range("E17").Value = '  -2.05%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.896.40'
$ws.Range("E18").Value = '  +1.23%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '8.86'
$ws.Range("E19").Value = '  -5.18%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '12.57'
$ws.Range("E20").Value = '  -0.05%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '376.40'
$ws.Range("E21").Value = '  -1.82%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '2.29'
$ws.Range("E22").Value = '  -1.09%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '4.15'
$ws.Range("E23").Value = '  -0.70%  '
$ws.Range("B24").Value = 'Dai'
$ws.Range("C24").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.999'
$ws.Range("E24").Value = '  -0.16%  '
$ws.Range("B25").Value = 'Litecoin'
$ws.Range("C25").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '71.09'
$ws.Range("E25").Value = '  -0.56%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '3.042.83'
$ws.Range("E26").Value = '  +1.31%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '4.18'
$ws.Range("E27").Value = '  -2.80%  '
$ws.Range("E28").Value = '  -2.13%  '
$ws.Range("E29").Value = '  +1.38%  '
$ws.Range("E30").Value = '  +0.06%  '
$ws.Range("E31").Value = '  -3.31%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '501.14'
$ws.Range("E32").Value = '  -8.05%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '7.68'
$ws.Range("E33").Value = '  -3.88%  '
$ws.Range("E34").Value = '  -3.00%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.999'
$ws.Range("E35").Value = '  +0.01%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '163.16'
$ws.Range("E36").Value = '  +0.78%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '20.01'
$ws.Range("E37").Value = '  -2.00%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '19.71'
$ws.Range("E38").Value = '  +2.06%  '
$ws.Range("E39").Value = '  -7.43%  '
$ws.Range("E40").Value = '  -0.15%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '179.62'
$ws.Range("E41").Value = '  -3.54%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.341'
$ws.Range("E42").Value = '  -1.54%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '4.99'
$ws.Range("E43").Value = '  -4.31%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.66'
$ws.Range("E44").Value = '  -3.71%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0903'
$ws.Range("E45").Value = '  +4.81%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.21'
$ws.Range("E46").Value = '  -5.53%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '40.03'
$ws.Range("E47").Value = '  +0.08%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.34'
$ws.Range("E48").Value = '  -3.29%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.574'
$ws.Range("E49").Value = '  -1.11%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '3.71'
$ws.Range("E50").Value = '  -2.58%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.651'
$ws.Range("E51").Value = '  +5.41%  '
